$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Test Environment description text in B11
$ws.Range("B11").Value = "OS: Windows `nBrowser: Chrome `n Network: Wi-Fi"

# Make column A (rows 2-15) bold
$ws.Range("A2:A15").Font.Bold = $true

# Set the active cell/selection to B11
$ws.Range("B11").Select()
